# Scale the "value" column (D) by 10,000 for all data rows (2-33).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $current = $cell.Value2
    if ($current -ne $null) {
        $scaled = [Math]::Round($current * 10000, 4)
        $cell.Value2 = $scaled
    }
}
